$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the values of E11 and H11
$e11 = $ws.Range("E11").Value2
$h11 = $ws.Range("H11").Value2
$ws.Range("E11").Value = $h11
$ws.Range("H11").Value = $e11

# Swap the values of E12 and H12
$e12 = $ws.Range("E12").Value2
$h12 = $ws.Range("H12").Value2
$ws.Range("E12").Value = $h12
$ws.Range("H12").Value = $e12

# Update the selected/active cell from E13 to C8
$ws.Range("C8").Select()
